$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update OBSERVACIONES (column O) for rows 15-29 with new commentary text ---
# (these become new shared-string entries, appended to sharedStrings.xml)
$ws.Range("O15").Value = 'El diseño de la BD debe ser revisado por el Líder Técnico para asegurar la escalabilidad (Req. No Funcional)'
$ws.Range("O16").Value = 'Tarea de backend. La dependencia con A23 es crítica para el inicio'
$ws.Range("O17").Value = 'Tarea clave para mitigar el Riesgo 7 (Disponibilidad). Se debe incluir la configuración de copias de seguridad'
$ws.Range("O18").Value = 'Punto de integración que suele generar errores técnicos. Se debe asignar tiempo para la corrección'
$ws.Range("O19").Value = 'Tarea crítica para el entregable principal (Mapa Emocional). Se debe considerar el Riesgo 8 (Rendimiento)'
$ws.Range("O20").Value = 'Primera validación end-to-end antes de las pruebas formales'
$ws.Range("O21").Value = 'Tarea clave para asegurar el cumplimiento de los requisitos funcionales'
$ws.Range("O22").Value = 'Mitiga el Riesgo 14 (Usabilidad Baja). Se debe coordinar el reclutamiento de participantes con anticipación'
$ws.Range("O23").Value = 'Tarea buffer para corregir los hallazgos de A29 y A30. La duración es corta; podría requerir más tiempo'
$ws.Range("O24").Value = 'La redacción debe ser simultánea al desarrollo para ser precisa y mitigar el Riesgo 10 (Mantenibilidad)'
$ws.Range("O25").Value = 'Esfuerzo necesario para asegurar el conocimiento del sistema'
$ws.Range("O26").Value = 'Tarea final de infraestructura. Se requiere la aprobación del Gerente de proyecto para el lanzamiento'
$ws.Range("O27").Value = 'Cierre formal con el Patrocinador. Se debe obtener la firma de aceptación del producto'
$ws.Range("O28").Value = 'Tarea clave para la entrega del valor. Se debe asegurar la alineación con la Investigación (A08)'
$ws.Range("O29").Value = 'Tarea final con el fin de estructurar y detallar lo realizado en el proyecto'

# --- Mark remaining in-flight activities (rows 19-29) as fully completed ---
$ws.Range("P19").Value = 100
$ws.Range("P20").Value = 100
$ws.Range("P21:P29").Value = 100

# --- Restore the working selection/view to match where the edits were made ---
$ws.Range("O28").Select()
